$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 30
$ws.Range("D3").Value = 30
$ws.Range("D4").Select()
